$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 gains new entries for the "Door" block type and the "Player" letter
$ws.Range("A17").Value = "Door"
$ws.Range("E17").Value = "Player"
$ws.Range("F17").Value = "U"

# New row 20 for the "Exit to next map" block type
$ws.Range("A20").Value = "Exit to next map"
$ws.Range("B20").Value = "x"

# Match the author's final selection
$ws.Range("B20").Select()
